# Weekly refresh of the "Fruta, Feria Lagunitas de Puerto Montt - Cereza" sheet.
#
# Effect of this commit: a new day's record is inserted at the top of the
# data block (row 17) and every existing record in rows 17..39 shifts down
# one row (into rows 18..40), which pushes the previous row 40 off the
# bottom of the block. Rows above the block (1..16) and below it (41..44)
# are untouched, as are the metadata columns A,B,C,E,F,G,H,I,J (constant
# for every row in this block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per record: D Fecha, K Variedad, L Calidad, M Volumen,
# N Precio minimo, O Precio maximo, P Precio promedio ponderado,
# Q Unidad de comercializacion, R Origen, S Precio $/Kg, T Kg/unidad.
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)

$firstRow = 17
$lastRow  = 39

# Snapshot the pre-edit values for the block that is about to shift, before
# any of it gets overwritten.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# Shift every snapshotted row down by one (process bottom-up so we never
# clobber a source row before it has been read).
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $dest = $r + 1
    foreach ($c in $cols) {
        $ws.Cells.Item($dest, $c).Value = $snapshot[$r][$c]
    }
}

# Write the brand-new record into row 17.
$ws.Cells.Item(17, 4).Value  = 44540                      # Fecha
$ws.Cells.Item(17, 11).Value = 'Santina'                  # Variedad
$ws.Cells.Item(17, 12).Value = 'Primera'                  # Calidad
$ws.Cells.Item(17, 13).Value = 800                        # Volumen
$ws.Cells.Item(17, 14).Value = 9000                       # Precio minimo
$ws.Cells.Item(17, 15).Value = 9000                       # Precio maximo
$ws.Cells.Item(17, 16).Value = 9000                       # Precio promedio ponderado
$ws.Cells.Item(17, 17).Value = '$/caja 8 kilos'           # Unidad de comercializacion
$ws.Cells.Item(17, 18).Value = 'Provincia de Curicó'      # Origen
$ws.Cells.Item(17, 19).Value = 1125                       # Precio $/Kg
$ws.Cells.Item(17, 20).Value = 8                          # Kg / unidad
